# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strikeouts) values replacing the old Strike# values in column G,
# row by row (row 1 is the header).
$newK = @{
    2  = 4
    3  = 1
    4  = 3
    5  = 0
    6  = 1
    7  = 3
    8  = 0
    10 = 3
    11 = 0
    12 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 0
    19 = 2
    20 = 3
    21 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
